$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet2 = $wb.Worksheets.Item("演出")
$sheet3 = $wb.Worksheets.Item("本地生活")
$sheet4 = $wb.Worksheets.Item("全部类型")

# --- Sheet: 展览 ---
$sheet1.Range("F3").Value = 3520
$sheet1.Range("F5").Value = 8234
$sheet1.Range("F7").Value = 92
$sheet1.Range("F8").Value = 2179
$sheet1.Range("F11").Value = 21
$sheet1.Range("F12").Value = 1198
$sheet1.Range("F13").Value = 59
$sheet1.Range("F15").Value = 20
$sheet1.Range("F16").Value = 585
$sheet1.Range("F18").Value = 3819
$sheet1.Range("F20").Value = 7283
$sheet1.Range("F22").Value = 55921
$sheet1.Range("F23").Value = 55921
$sheet1.Range("F24").Value = 4478
$sheet1.Range("F26").Value = 1040
$sheet1.Range("F27").Value = 873
$sheet1.Range("F28").Value = 433
$sheet1.Range("F29").Value = 91
$sheet1.Range("F30").Value = 889
$sheet1.Range("F31").Value = 598
$sheet1.Range("F32").Value = 3328
$sheet1.Range("F33").Value = 584
$sheet1.Range("F34").Value = 45
$sheet1.Range("F35").Value = 28
$sheet1.Range("F37").Value = 1218
$sheet1.Range("F38").Value = 1184
$sheet1.Range("F39").Value = 158
$sheet1.Range("F40").Value = 192
$sheet1.Range("F41").Value = 1070
$sheet1.Range("F42").Value = 705
$sheet1.Range("F43").Value = 7
$sheet1.Range("F44").Value = 767
$sheet1.Range("F45").Value = 166
$sheet1.Range("F47").Value = 167

# --- Sheet: 演出 ---
$sheet2.Range("F8").Value = 174
$sheet2.Range("F10").Value = 51
$sheet2.Range("F11").Value = 47
$sheet2.Range("F12").Value = 112
$sheet2.Range("F14").Value = 40
$sheet2.Range("F15").Value = 172
$sheet2.Range("F16").Value = 7480
$sheet2.Range("F35").Value = 31
$sheet2.Range("F44").Value = 65
$sheet2.Range("F45").Value = 29

# --- Sheet: 本地生活 ---
$sheet3.Range("F4").Value = 2291
$sheet3.Range("F5").Value = 1552
$sheet3.Range("F8").Value = 2336
$sheet3.Range("F10").Value = 1671
$sheet3.Range("F11").Value = 161
$sheet3.Range("F12").Value = 86
$sheet3.Range("F15").Value = 168

# --- Sheet: 全部类型 ---
$sheet4.Range("F3").Value = 3520
$sheet4.Range("F5").Value = 8234
$sheet4.Range("F6").Value = 1552
$sheet4.Range("F8").Value = 2336
$sheet4.Range("F9").Value = 1671
$sheet4.Range("F10").Value = 161
$sheet4.Range("F11").Value = 86
$sheet4.Range("F13").Value = 92
$sheet4.Range("F16").Value = 59
$sheet4.Range("F17").Value = 20
$sheet4.Range("F18").Value = 585
$sheet4.Range("F20").Value = 7283
$sheet4.Range("F21").Value = 55921
$sheet4.Range("F23").Value = 51
$sheet4.Range("F24").Value = 4478
$sheet4.Range("F25").Value = 1040
$sheet4.Range("F26").Value = 433
$sheet4.Range("F27").Value = 91
$sheet4.Range("F28").Value = 598
$sheet4.Range("F29").Value = 112
$sheet4.Range("F30").Value = 3328
$sheet4.Range("F31").Value = 584
$sheet4.Range("F32").Value = 45
$sheet4.Range("F33").Value = 28
$sheet4.Range("F35").Value = 1218
$sheet4.Range("F37").Value = 158
$sheet4.Range("F38").Value = 192
$sheet4.Range("F39").Value = 1070
$sheet4.Range("F40").Value = 705
$sheet4.Range("F41").Value = 767
$sheet4.Range("F42").Value = 166
$sheet4.Range("F44").Value = 167
$sheet4.Range("F47").Value = 31
